# WatchlistTestData.xlsx edit script
# - Rewrites the DESCRIPTION (col B) text for every test row with clearer
#   "Verify that ..." wording.
# - Populates / extends the VALIDATIONS (col J) values.
# - Clears the old STATUS (col L) column (the PASS markers go away).
# - Turns on word-wrap for columns B, H, J (and the B1/J1 header cells)
#   and enlarges several row heights so the wrapped text is visible.
# - Minor column width tweaks for H and J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values
# ---------------------------------------------------------------------

# Row 2 - S1_TC_T1 (search)
$ws.Range("B2").Value = "Verify that to get document for adding into users watchlist"
$ws.Range("J2").Value = "status=200"
$ws.Range("L2").ClearContents()

# Row 3 - S1_TC_T2 (get watchlist)
$ws.Range("B3").Value = "Verify that to get user watchlist "
$ws.Range("J3").Value = "status=200"
$ws.Range("L3").ClearContents()

# Row 4 - S1_TC_T3 (add item to watchlist)
$ws.Range("B4").Value = "Verify that to add new item in to users watch watchlist"
$ws.Range("L4").ClearContents()

# Row 5 - S1_TC_T4 (get watchlist again / validate added item)
$ws.Range("B5").Value = "Verify that to get user watchlist and validate user added article should be avilable in his watch list"
$ws.Range("J5").Value = "status=200||userId=(SYS_USER1)||items.itemId=(S1_TC_T1_hits.hits._id)"
$ws.Range("L5").ClearContents()

# Row 6 - S1_TC_T5 (who has item)
$ws.Range("B6").Value = "Verify that to get all users watching an article are returned based on article id"
$ws.Range("J6").Value = "status=200||userId=(SYS_USER1)"
$ws.Range("L6").ClearContents()

# Row 7 - S1_TC_T6 (delete item)
$ws.Range("B7").Value = "Verify that user is able to delete a document from watchlist"
$ws.Range("L7").ClearContents()

# Row 8 - S1_TC_T7 (get watchlist / validate deleted item)
$ws.Range("B8").Value = "Verify that to get user watchlist and validate user deleted article should not be avilable in his watch list"
$ws.Range("J8").Value = "status=200"
$ws.Range("L8").ClearContents()

# ---------------------------------------------------------------------
# 2. Word wrap formatting for columns B, H, J (header + data rows)
# ---------------------------------------------------------------------

$ws.Range("B1:B8").WrapText = $true
$ws.Range("H1:H8").WrapText = $true
$ws.Range("J1:J8").WrapText = $true

# ---------------------------------------------------------------------
# 3. Row heights for the rows whose description text now wraps to
#    multiple lines.
# ---------------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 45

# ---------------------------------------------------------------------
# 4. Column width tweaks (H gets narrower, J gets wider now that it
#    wraps). Values are expressed in "character" units as used by the
#    ColumnWidth COM property.
# ---------------------------------------------------------------------

$ws.Columns.Item(8).ColumnWidth = 6.333333333333333
$ws.Columns.Item(10).ColumnWidth = 24.333333333333332

Write-Host "WatchlistTestData edits applied"
